$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.291
$ws.Range("D4").Value = 0.227
$ws.Range("E4").Value = 0.149
$ws.Range("H4").Value = 0.187
$ws.Range("J4").Value = 0.109
$ws.Range("K4").Value = 0.351
$ws.Range("L4").Value = 0.102
$ws.Range("M4").Value = 0.319
$ws.Range("N4").Value = 0.279
$ws.Range("P4").Value = 0.137
$ws.Range("Q4").Value = 0.527
$ws.Range("R4").Value = 0.215
$ws.Range("S4").Value = 0.464
$ws.Range("T4").Value = 0.291
$ws.Range("W4").Value = 0.242
$ws.Range("Y4").Value = 0.21
$ws.Range("Z4").Value = 0.436
$ws.Range("AA4").Value = 0.128
$ws.Range("AB4").Value = 0.357
$ws.Range("AE4").Value = 0.076
$ws.Range("AF4").Value = 0.73
$ws.Range("AH4").Value = 0.308
$ws.Range("AI4").Value = 0.676
$ws.Range("AJ4").Value = 0.165
$ws.Range("AK4").Value = 0.406
$ws.Range("AL4").Value = 0.708
$ws.Range("AM4").Value = 0.117
$ws.Range("AN4").Value = 0.342
$ws.Range("AO4").Value = 0.705
$ws.Range("B5").Value = 0.8110000000000001
$ws.Range("C5").Value = 0.153
$ws.Range("D5").Value = 0.392
$ws.Range("E5").Value = 0.649
$ws.Range("F5").Value = 0.228
$ws.Range("G5").Value = 0.477
$ws.Range("H5").Value = 0.8110000000000001
$ws.Range("I5").Value = 0.153
$ws.Range("J5").Value = 0.392
$ws.Range("K5").Value = 0.676
$ws.Range("L5").Value = 0.219
$ws.Range("M5").Value = 0.468
$ws.Range("N5").Value = 0.865
$ws.Range("O5").Value = 0.117
$ws.Range("P5").Value = 0.342
$ws.Range("Q5").Value = 0.595
$ws.Range("R5").Value = 0.241
$ws.Range("S5").Value = 0.491
$ws.Range("T5").Value = 0.595
$ws.Range("U5").Value = 0.241
$ws.Range("V5").Value = 0.491
$ws.Range("W5").Value = 0.73
$ws.Range("X5").Value = 0.197
$ws.Range("Y5").Value = 0.444
$ws.Range("Z5").Value = 0.8110000000000001
$ws.Range("AA5").Value = 0.153
$ws.Range("AB5").Value = 0.392
$ws.Range("AC5").Value = 0.784
$ws.Range("AD5").Value = 0.169
$ws.Range("AE5").Value = 0.412
$ws.Range("AF5").Value = 0.973
$ws.Range("AG5").Value = 0.026
$ws.Range("AH5").Value = 0.162
$ws.Range("AI5").Value = 0.784
$ws.Range("AJ5").Value = 0.169
$ws.Range("AK5").Value = 0.412
$ws.Range("AL5").Value = 0.919
$ws.Range("AM5").Value = 0.075
$ws.Range("AN5").Value = 0.273
$ws.Range("AO5").Value = 0.892
$ws.Range("B6").Value = 0.428
$ws.Range("E6").Value = 0.242
$ws.Range("H6").Value = 0.304
$ws.Range("K6").Value = 0.462
$ws.Range("N6").Value = 0.422
$ws.Range("Q6").Value = 0.5590000000000001
$ws.Range("T6").Value = 0.391
$ws.Range("W6").Value = 0.363
$ws.Range("Z6").Value = 0.5669999999999999
$ws.Range("AC6").Value = 0.222
$ws.Range("AF6").Value = 0.834
$ws.Range("AI6").Value = 0.726
$ws.Range("AL6").Value = 0.8
$ws.Range("AO6").Value = 0.787
$ws.Range("B7").Value = 0.597
$ws.Range("E7").Value = 0.388
$ws.Range("H7").Value = 0.486
$ws.Range("K7").Value = 0.57
$ws.Range("N7").Value = 0.609
$ws.Range("Q7").Value = 0.58
$ws.Range("T7").Value = 0.492
$ws.Range("W7").Value = 0.52
$ws.Range("Z7").Value = 0.6919999999999999
$ws.Range("AC7").Value = 0.389
$ws.Range("AF7").Value = 0.912
$ws.Range("AI7").Value = 0.76
$ws.Range("AL7").Value = 0.867
$ws.Range("AO7").Value = 0.846
$ws.Range("B8").Value = 0.754
$ws.Range("C8").Value = 0.153
$ws.Range("D8").Value = 0.391
$ws.Range("E8").Value = 0.546
$ws.Range("F8").Value = 0.196
$ws.Range("H8").Value = 0.705
$ws.Range("I8").Value = 0.156
$ws.Range("J8").Value = 0.395
$ws.Range("K8").Value = 0.602
$ws.Range("L8").Value = 0.2
$ws.Range("M8").Value = 0.447
$ws.Range("N8").Value = 0.769
$ws.Range("O8").Value = 0.126
$ws.Range("P8").Value = 0.354
$ws.Range("Q8").Value = 0.5649999999999999
$ws.Range("R8").Value = 0.227
$ws.Range("S8").Value = 0.476
$ws.Range("T8").Value = 0.514
$ws.Range("U8").Value = 0.204
$ws.Range("V8").Value = 0.452
$ws.Range("W8").Value = 0.653
$ws.Range("X8").Value = 0.182
$ws.Range("Y8").Value = 0.426
$ws.Range("Z8").Value = 0.744
$ws.Range("AA8").Value = 0.152
$ws.Range("AB8").Value = 0.39
$ws.Range("AC8").Value = 0.672
$ws.Range("AD8").Value = 0.169
$ws.Range("AE8").Value = 0.412
$ws.Range("AF8").Value = 0.89
$ws.Range("AG8").Value = 0.047
$ws.Range("AH8").Value = 0.218
$ws.Range("AI8").Value = 0.774
$ws.Range("AJ8").Value = 0.169
$ws.Range("AK8").Value = 0.411
$ws.Range("AL8").Value = 0.889
$ws.Range("AM8").Value = 0.08
$ws.Range("AN8").Value = 0.282
$ws.Range("AO8").Value = 0.851
$ws.Range("B9").Value = 0.676
$ws.Range("C9").Value = 0.219
$ws.Range("D9").Value = 0.468
$ws.Range("E9").Value = 0.432
$ws.Range("F9").Value = 0.245
$ws.Range("G9").Value = 0.495
$ws.Range("H9").Value = 0.595
$ws.Range("I9").Value = 0.241
$ws.Range("J9").Value = 0.491
$ws.Range("K9").Value = 0.514
$ws.Range("N9").Value = 0.649
$ws.Range("O9").Value = 0.228
$ws.Range("P9").Value = 0.477
$ws.Range("Q9").Value = 0.514
$ws.Range("R9").Value = 0.25
$ws.Range("S9").Value = 0.5
$ws.Range("T9").Value = 0.405
$ws.Range("U9").Value = 0.241
$ws.Range("V9").Value = 0.491
$ws.Range("W9").Value = 0.541
$ws.Range("X9").Value = 0.248
$ws.Range("Y9").Value = 0.498
$ws.Range("Z9").Value = 0.649
$ws.Range("AA9").Value = 0.228
$ws.Range("AB9").Value = 0.477
$ws.Range("AC9").Value = 0.5679999999999999
$ws.Range("AD9").Value = 0.245
$ws.Range("AE9").Value = 0.495
$ws.Range("AF9").Value = 0.757
$ws.Range("AG9").Value = 0.184
$ws.Range("AH9").Value = 0.429
$ws.Range("AI9").Value = 0.757
$ws.Range("AJ9").Value = 0.184
$ws.Range("AK9").Value = 0.429
$ws.Range("AL9").Value = 0.838
$ws.Range("AM9").Value = 0.136
$ws.Range("AN9").Value = 0.369
$ws.Range("AO9").Value = 0.784
$ws.Range("B10").Value = 0.8110000000000001
$ws.Range("C10").Value = 0.153
$ws.Range("D10").Value = 0.392
$ws.Range("E10").Value = 0.5679999999999999
$ws.Range("F10").Value = 0.245
$ws.Range("G10").Value = 0.495
$ws.Range("H10").Value = 0.73
$ws.Range("I10").Value = 0.197
$ws.Range("J10").Value = 0.444
$ws.Range("K10").Value = 0.676
$ws.Range("L10").Value = 0.219
$ws.Range("M10").Value = 0.468
$ws.Range("N10").Value = 0.838
$ws.Range("O10").Value = 0.136
$ws.Range("P10").Value = 0.369
$ws.Range("Q10").Value = 0.595
$ws.Range("R10").Value = 0.241
$ws.Range("S10").Value = 0.491
$ws.Range("T10").Value = 0.595
$ws.Range("U10").Value = 0.241
$ws.Range("V10").Value = 0.491
$ws.Range("W10").Value = 0.73
$ws.Range("X10").Value = 0.197
$ws.Range("Y10").Value = 0.444
$ws.Range("Z10").Value = 0.8110000000000001
$ws.Range("AA10").Value = 0.153
$ws.Range("AB10").Value = 0.392
$ws.Range("AC10").Value = 0.676
$ws.Range("AD10").Value = 0.219
$ws.Range("AE10").Value = 0.468
$ws.Range("AF10").Value = 0.973
$ws.Range("AG10").Value = 0.026
$ws.Range("AH10").Value = 0.162
$ws.Range("AI10").Value = 0.784
$ws.Range("AJ10").Value = 0.169
$ws.Range("AK10").Value = 0.412
$ws.Range("AL10").Value = 0.919
$ws.Range("AM10").Value = 0.075
$ws.Range("AN10").Value = 0.273
$ws.Range("AO10").Value = 0.892
$ws.Range("B11").Value = 0.8110000000000001
$ws.Range("C11").Value = 0.153
$ws.Range("D11").Value = 0.392
$ws.Range("E11").Value = 0.649
$ws.Range("F11").Value = 0.228
$ws.Range("G11").Value = 0.477
$ws.Range("H11").Value = 0.8110000000000001
$ws.Range("I11").Value = 0.153
$ws.Range("J11").Value = 0.392
$ws.Range("K11").Value = 0.676
$ws.Range("L11").Value = 0.219
$ws.Range("M11").Value = 0.468
$ws.Range("N11").Value = 0.865
$ws.Range("O11").Value = 0.117
$ws.Range("P11").Value = 0.342
$ws.Range("Q11").Value = 0.595
$ws.Range("R11").Value = 0.241
$ws.Range("S11").Value = 0.491
$ws.Range("T11").Value = 0.595
$ws.Range("U11").Value = 0.241
$ws.Range("V11").Value = 0.491
$ws.Range("W11").Value = 0.73
$ws.Range("X11").Value = 0.197
$ws.Range("Y11").Value = 0.444
$ws.Range("Z11").Value = 0.8110000000000001
$ws.Range("AA11").Value = 0.153
$ws.Range("AB11").Value = 0.392
$ws.Range("AC11").Value = 0.73
$ws.Range("AD11").Value = 0.197
$ws.Range("AE11").Value = 0.444
$ws.Range("AF11").Value = 0.973
$ws.Range("AG11").Value = 0.026
$ws.Range("AH11").Value = 0.162
$ws.Range("AI11").Value = 0.784
$ws.Range("AJ11").Value = 0.169
$ws.Range("AK11").Value = 0.412
$ws.Range("AL11").Value = 0.919
$ws.Range("AM11").Value = 0.075
$ws.Range("AN11").Value = 0.273
$ws.Range("AO11").Value = 0.892
$ws.Range("B12").Value = 1.233
$ws.Range("C12").Value = 0.312
$ws.Range("D12").Value = 0.5590000000000001
$ws.Range("E12").Value = 1.667
$ws.Range("F12").Value = 1.139
$ws.Range("G12").Value = 1.067
$ws.Range("H12").Value = 1.6
$ws.Range("I12").Value = 1.373
$ws.Range("J12").Value = 1.172
$ws.Range("K12").Value = 1.4
$ws.Range("L12").Value = 0.5600000000000001
$ws.Range("M12").Value = 0.748
$ws.Range("N12").Value = 1.406
$ws.Range("O12").Value = 0.616
$ws.Range("P12").Value = 0.785
$ws.Range("Z12").Value = 1.267
$ws.Range("AA12").Value = 0.329
$ws.Range("AB12").Value = 0.573
$ws.Range("AC12").Value = 1.793
$ws.Range("AD12").Value = 2.44
$ws.Range("AE12").Value = 1.562
$ws.Range("AF12").Value = 1.25
$ws.Range("AG12").Value = 0.243
$ws.Range("AH12").Value = 0.493
$ws.Range("AI12").Value = 1.034
$ws.Range("AJ12").Value = 0.033
$ws.Range("AK12").Value = 0.182
$ws.Range("AL12").Value = 1.088
$ws.Range("AM12").Value = 0.08
$ws.Range("AN12").Value = 0.284
$ws.Range("AO12").Value = 1.124
$ws.Range("B13").Value = 3.514
$ws.Range("C13").Value = 1.439
$ws.Range("D13").Value = 1.2
$ws.Range("E13").Value = 4.613
$ws.Range("F13").Value = 0.431
$ws.Range("G13").Value = 0.656
$ws.Range("H13").Value = 4.629
$ws.Range("I13").Value = 0.633
$ws.Range("J13").Value = 0.796
$ws.Range("K13").Value = 2.265
$ws.Range("L13").Value = 0.606
$ws.Range("M13").Value = 0.779
$ws.Range("N13").Value = 3.27
$ws.Range("O13").Value = 0.738
$ws.Range("P13").Value = 0.859
$ws.Range("Z13").Value = 2.559
$ws.Range("AA13").Value = 2.894
$ws.Range("AB13").Value = 1.701
$ws.Range("AC13").Value = 6.378
$ws.Range("AD13").Value = 2.181
$ws.Range("AE13").Value = 1.477
$ws.Range("AF13").Value = 1.622
$ws.Range("AG13").Value = 0.722
$ws.Range("AH13").Value = 0.85
$ws.Range("AI13").Value = 1.297
$ws.Range("AJ13").Value = 0.371
$ws.Range("AK13").Value = 0.609
$ws.Range("AL13").Value = 1.568
$ws.Range("AM13").Value = 0.732
$ws.Range("AN13").Value = 0.856
$ws.Range("AO13").Value = 1.496
